# Logs and connection feedback adjustment
# Insert two new journal rows before the "43614" date row, fill them in,
# and move the selection to the second new cell - matching how the author
# would have entered new log lines in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 59 (the "43614" entry), pushing the
# existing tail of the table (dates + totals) down by two rows.
$ws.Rows.Item(59).Resize(2).Insert()

# Excel's row insert doesn't always carry the surrounding formatting
# through this host, so explicitly pull the format from the row above
# (same look as every other data row in the table).
$ws.Range("A58:D58").Copy()
$ws.Range("A59:D60").PasteSpecial(-4122)

# New row 59: Debug entry, 1.5h, with a remark.
$ws.Cells.Item(59, 2).Value = "Debug"
$ws.Cells.Item(59, 3).Value = 1.5
$ws.Cells.Item(59, 4).Value = "Première version fonctionnelle"

# New row 60: another entry.
$ws.Cells.Item(60, 2).Value = "Test"

# The total row's SUM formula auto-expands to include the two inserted
# rows (same as Excel would do), so no manual formula rewrite is needed.

$ws.Range("B60").Select()
